$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph, with "Meta description" bold and the rest
#    of the sentence in normal formatting.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaStart = $metaPara.Range.Start
$insertPoint = $d.Range($metaStart, $metaStart)
$insertPoint.InsertAfter("Meta description: Immersive gameplay with a Western twist. Review covers symbols, paylines, RTP, and features. Play Crazy Colt free.")

# Re-resolve the paragraph start (it may have shifted) and bold the
# "Meta description" label only.
$metaStart = $d.Paragraphs.Item(2).Range.Start
$labelRange = $d.Range($metaStart, $metaStart + 16)
$labelRange.Bold = 1

# Match the document's convention of a leading empty run in body
# paragraphs by inserting an empty run at the very start.
$leadRange = $d.Range($metaStart, $metaStart)
$leadRange.InsertBefore("")

# ---------------------------------------------------------------------
# 2) Near the end of the document: drop the duplicated bold title
#    paragraph, and replace the italic "meta description" paragraph's
#    text with the new feature-image prompt (keeping italic formatting).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
if ($dupTitlePara.Range.Text.TrimEnd() -ne "Play Crazy Colt Free - Review of JVC's Western-Themed Slot") {
    throw "Unexpected paragraph content where the duplicated title was expected: $($dupTitlePara.Range.Text)"
}
$dupTitleRange = $d.Range($dupTitlePara.Range.Start, $dupTitlePara.Range.End)
$dupTitleRange.Delete()

$count = $d.Paragraphs.Count
$imgPromptPara = $d.Paragraphs.Item($count)
if ($imgPromptPara.Range.Text.TrimEnd() -ne "Immersive gameplay with a Western twist. Review covers symbols, paylines, RTP, and features. Play Crazy Colt free.") {
    throw "Unexpected paragraph content where the italic meta paragraph was expected: $($imgPromptPara.Range.Text)"
}
$imgPromptRange = $d.Range($imgPromptPara.Range.Start, $imgPromptPara.Range.End)
$imgPromptRange.Text = "Create a feature image that stands out with a cartoon-style happy Maya warrior. The image should be set against the dusty landscape of Arizona, mirroring the theme of the game, ""Crazy Colt"". The warrior should be wearing glasses to highlight their intelligence, adding a unique character trait to the image. Make sure to convey a sense of excitement and adventure in the image, as the slot game is all about thrilling moments and big wins. Use bold colors and sharp lines to make the image stand out, attracting attention to this game and enticing players to take a spin."

Write-Host "Done."
